$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray duplicate hidden chart-data defined names (v1.2 / v1.3) -
# only v1.0 / v1.1 are actually used by the box-whisker chart.
$wb.Names.Item("_xlchart.v1.3").Delete()
$wb.Names.Item("_xlchart.v1.2").Delete()

# New "Mean increase" / "Median increase" headers (row 18, bold like the
# other D/E stat headers) with the increase-vs-baseline formulas beneath
# them (row 19) referencing the existing mean (E4) / median (E10) cells.
$ws.Range("D18").Value = "Mean increase"
$ws.Range("D18").Font.Bold = $true

$ws.Range("F18").Value = "Median increase"
$ws.Range("F18").Font.Bold = $true

$ws.Range("D19").Formula = "=((E4 / 114.202998) * 100) - 100"
$ws.Range("D19").ClearFormats()

$ws.Range("F19").Formula = "=((E10 / 113.658804) * 100) - 100"
$ws.Range("F19").ClearFormats()

# Move the selection/active cell the way the author left it.
$ws.Range("D20").Select() | Out-Null
